$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.283.55'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '2.995.85'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.55'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.23'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.26%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.428'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.09'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.367'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.63%  '
$ws.Range('D12').Value = '3.509.95'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.32'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('D16').Value = '56.206.34'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').Value = '2.994.95'
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.91'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.84'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('E20').Value = '  +2.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '330.27'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.495'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.41'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.07%  '
$ws.Range('D25').Value = '3.124.38'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.164'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').Value = '0.0₃0929'
$ws.Range('E28').Value = '  +5.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.31'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.84'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.78'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.27'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '152.35'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.70%  '
$ws.Range('E35').Value = '  -2.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.56'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +10.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.79'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.21'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.89%  '
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').Value = '3.025.05'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.35'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.76'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.653'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.30%  '
$ws.Range('D45').Value = '2.192.13'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0238'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.80'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.916'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.51'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0846'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.96%  '
